# Update the "Förändrad" (changed) date column (C) for all data rows
# on the "Avverkningsanmälningar" sheet from 2023-09-16 (45185) to
# 2023-10-05 (45204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
